$d = $word.ActiveDocument

# 1) "The program then automatically enters ..." -> reworded with menu reprint + "choice = -1" phrasing
$d.Content.Find.Execute(
    "The program then automatically enters " + [char]0x201C + "-1" + [char]0x201D + " as the input.  The program again displays " + [char]0x201C + "wrong answer" + [char]0x201D + " 500 times again.  On the fourth iteration, the program prints",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The program then reprints the menu options, and then automatically inputs " + [char]0x201C + "-1" + [char]0x201D + " while printing " + [char]0x201C + "choice = -1" + [char]0x201D + ".  The program again displays " + [char]0x201C + "wrong answer" + [char]0x201D + " 500 times.  The process repeats until in the fourth iteration, the program prints",
    2
)

# 2) " again" 1,000 times to the screen. After the program automatically inputs ..." -> reworded tail
$d.Content.Find.Execute(
    " again" + [char]0x201D + " 1,000 times to the screen. After the program automatically inputs " + [char]0x201C + "-1" + [char]0x201D + " once more, the words " + [char]0x201C + "I warned you" + [char]0x201D + " are output multiple times",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " again" + [char]0x201D + " 1,000 times to the screen after reprinting the menu. The program automatically inputs " + [char]0x201C + "-1" + [char]0x201D + " once more while printing " + [char]0x201C + "choice = -1" + [char]0x201D + ", and again outputs " + [char]0x201C + "wrong answer" + [char]0x201D + " 500 times. On the next iteration, after printing the menu, the words " + [char]0x201C + "I warned you" + [char]0x201D + " are output 1000 times",
    2
)

# 3) "that loops" -> "that loops for a 1000 times"
$d.Content.Find.Execute(
    "that loops",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "that loops for a 1000 times",
    2
)
